$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# ---------------------------------------------------------------------------
# Weekly Quantity: append two new weekly rows (rows 3 and 4)
# ---------------------------------------------------------------------------
$wsWeekly.Range("A2:B2").Copy($wsWeekly.Range("A3:B4"))
$wsWeekly.Range("A3").Value = 45662.99999999999
$wsWeekly.Range("B3").Value = 8
$wsWeekly.Range("A4").Value = 45676.99999999999
$wsWeekly.Range("B4").Value = 5

# ---------------------------------------------------------------------------
# Monthly Trend: append one new monthly row (row 3)
# ---------------------------------------------------------------------------
$wsMonthly.Range("A2:B2").Copy($wsMonthly.Range("A3:B3"))
$wsMonthly.Range("A3").Value = 45688.99999999999
$wsMonthly.Range("B3").Value = 13

# ---------------------------------------------------------------------------
# PO Forecast: brand new sheet placed after "Monthly Trend"
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Stamp formatting: header row style from Weekly Quantity's header,
# and the date-value row style across all 11 data rows.
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("A1:B1"))
$wsWeekly.Range("A2:B2").Copy($wsForecast.Range("A2:B12"))

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"

$forecastDates = @(45613.99999999999, 45662.99999999999, 45676.99999999999, 45683.99999999999, 45690.99999999999, 45697.99999999999, 45704.99999999999, 45711.99999999999, 45718.99999999999, 45725.99999999999, 45732.99999999999)
$forecastValues = @(37, 11, 3, 0, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $forecastDates.Length; $i++) {
    $row = $i + 2
    $wsForecast.Range("A$row").Value = $forecastDates[$i]
    $wsForecast.Range("B$row").Value = $forecastValues[$i]
}

# Restore the originally active sheet/selection.
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select() | Out-Null

